# Adiciona atributos e métodos na definição UML da classe "Escopo":
#   - String[] ?? comandos
#   - ArrayList variaveis
#   + void processa()

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "String[] ??"
$ws.Range("F4").Value = "comandos"

$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "ArrayList"
$ws.Range("F5").Value = "variaveis"

$ws.Range("D7").Value = "+"
$ws.Range("E7").Value = "void"
$ws.Range("F7").Value = "processa()"
